$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (KILLS) - convert text-number cells to real numbers
$aValues = @{
    2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0;
    11=1; 12=1;
    13=2; 14=2; 15=2; 16=2; 17=2; 18=2; 19=2; 20=2; 21=2; 22=2; 23=2;
    24=3;
    25=4; 26=4; 27=4; 28=4; 29=4; 30=4; 31=4; 32=4; 33=4; 34=4; 35=4;
    36=5; 37=5; 38=5; 39=5; 40=5; 41=5
}

foreach ($row in $aValues.Keys) {
    $ws.Range("A$row").Value = $aValues[$row]
}

# Column E (DEATHS) - convert text-number cells to real numbers (also fixes E31 typo 12 -> 2)
$eValues = @{
    2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 13=0; 14=0; 15=0; 16=0; 17=0;
    18=1; 19=1; 20=1; 21=1; 22=1; 23=1; 24=1; 25=1; 26=1; 27=1; 28=1; 29=1;
    30=2; 31=2; 32=2; 33=2; 34=2; 35=2; 36=2;
    37=3; 38=3; 39=3; 40=3; 41=3
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}

# F41 (ASSISTS) text fix from "erro" to "6" (stays a text value, not a number)
$ws.Range("F41").Value = "'6"
